# Applies the edits described by the diff:
#  - Hyperparameter Optimization!Q24: "All" -> "All *"
#  - Hyperparameter Optimization! rows 25,44,68-75, cols S:V: fill in previously-empty values
#  - CNNLSTM SH! rows 4,8,13,14,15,19,20,21, cols B:F and H:L: fill in previously-empty values
#  - CNNLSTM SH!L31: fill in previously-empty value (dependent formulas in M31/L32/M32 recalc automatically)
#  - Selection / view changes on SARIMAX, CNNLSTM SH and Hyperparameter Optimization sheets

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. CNNLSTM SH sheet: fill in the B:F / H:L blocks and L31
# ---------------------------------------------------------------------------
$wsCnnLstmSh = $wb.Worksheets.Item("CNNLSTM SH")

$cnnLstmShData = @{
    4  = @{ B=0.22991296682034901; C=0.38238020028036601; D=0.316772729333612;  E=0.319113515541112;   F=0.242870078348318;
            H=0.46091266474450399; I=0.37971040075847901; J=0.48448361693066599; K=0.63138786230808797; L=0.471823824387955 }
    8  = @{ B=0.59539785574497595; C=0.39580970455407699; D=0.62325278302609399; E=0.575362929923661;   F=0.55851517942149898;
            H=0.25112878687831702; I=0.31776313049497101; J=0.31977581961944002; K=0.35589093623428097; L=0.25654343858890599 }
    13 = @{ B=0.30547930957932501; C=0.29671517141179199; D=0.38608664105482798; E=0.77493506811865498; F=0.41423679997721402;
            H=0.46480240450395899; I=0.402020327243088;    J=0.47191390734466299; K=0.70998165388802703; L=0.30803363635861603 }
    14 = @{ B=0.50146764320834303; C=0.15808596850718501; D=0.734342937725324;   E=0.90949135692244798; F=0.761584049059715;
            H=1.1752916757001199;  I=0.411144331442725;    J=1.0849470174759099;  K=0.78268124750582702; L=0.93033754021894099 }
    15 = @{ B=0.46633325098598899; C=0.53164927700737297; D=0.63850139288685603; E=1.7643581092013301;  F=0.30090417907529199;
            H=2.3433267063471099;  I=0.68894307405754396; J=1.1591996622109;     K=0.70974098038565403; L=0.715846561090692 }
    19 = @{ B=0.31347803287470799; C=0.38905770782642801; D=0.546419422863618;   E=0.54863573177676905; F=0.370153307112913;
            H=0.30298812620099602; I=0.38919164453733601; J=0.58889396493185897; K=0.590925664985952;   L=0.30848371481190701 }
    20 = @{ B=0.64934536471215898; C=0.493202713899487;   D=1.2365408059999401;  E=0.66716813448911505; F=1.1006613261270399;
            H=0.61887762881156505; I=0.72259952317502296; J=1.2523959819358199;  K=0.77086694931597;    L=0.91935349333691097 }
    21 = @{ B=1.39381378634936;    C=0.84542474534971801; D=1.30751034591801;    E=1.1314087062235001;  F=0.93741812459009799;
            H=1.2580063814756799;  I=1.3659570278842801;  J=1.133470445732;      K=0.63444760715533399; L=0.69492734663255495 }
}

foreach ($row in $cnnLstmShData.Keys) {
    $cells = $cnnLstmShData[$row]
    foreach ($col in $cells.Keys) {
        $wsCnnLstmSh.Range("$col$row").Value = $cells[$col]
    }
}

# Previously-empty L31 (dependent cells M31, L32, M32 hold formulas that will
# recompute automatically on recalculation).
$wsCnnLstmSh.Range("L31").Value = 1.1293309972672301

# ---------------------------------------------------------------------------
# 2. Hyperparameter Optimization sheet: Q24 text change + S:V fills
# ---------------------------------------------------------------------------
$wsHyperOpt = $wb.Worksheets.Item("Hyperparameter Optimization")

$wsHyperOpt.Range("Q24").Value = "All *"

$hyperOptData = @{
    25 = @{ S=5.6768941274403302;  T=5.1989356565323703;  U=1.88962965674212;     V=-0.83036371922891306 }
    44 = @{ S=1.88441647665034;    T=1.5500394817199501;  U=1.65194906068285;     V=0.0573353237453516 }
    68 = @{ S=0.29820989806475101; T=0.29820989806475101; U=0.065826197464093503; V=0.96005124343230597 }
    69 = @{ S=0.54966769053406095; T=0.45292373696240501; U=0.104614083056826;    V=0.85317659209757701 }
    70 = @{ S=0.30022042236318303; T=0.243343818239867;   U=0.0620133573197744;   V=0.95968074899584599 }
    71 = @{ S=0.485663673825938;   T=0.37491296203028002; U=0.087352969440431305; V=0.88910603037918701 }
    72 = @{ S=0.59627807698144497; T=0.52232116266375905; U=0.129592088257662;    V=0.76195480754216605 }
    73 = @{ S=0.79534921707419204; T=0.65246940823983002; U=0.136674488514721;    V=0.65384653490961697 }
    74 = @{ S=0.77009236672821302; T=0.61049838719765304; U=0.12272515504937399;  V=0.69460931318921404 }
    75 = @{ S=0.823880715051587;   T=0.66734719294219302; U=0.129396832789067;    V=0.45297592084250499 }
}

foreach ($row in $hyperOptData.Keys) {
    $cells = $hyperOptData[$row]
    foreach ($col in $cells.Keys) {
        $wsHyperOpt.Range("$col$row").Value = $cells[$col]
    }
}

# ---------------------------------------------------------------------------
# 3. Selections / active-sheet views.
#    Activate in this order so the workbook's final active sheet/tab matches
#    the target file (Hyperparameter Optimization stays active/tabSelected).
# ---------------------------------------------------------------------------
$wsSarimax = $wb.Worksheets.Item("SARIMAX")
$wsSarimax.Activate()
$wsSarimax.Range("H11").Select()

$wsCnnLstmSh.Activate()
$wsCnnLstmSh.Range("H5").Select()

$wsHyperOpt.Activate()
$wsHyperOpt.Range("S26").Select()
